$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.530.81'
$ws.Range("E2").Value = '  +2.35%  '

$ws.Range("D3").Value = '1.597.41'
$ws.Range("E3").Value = '  +1.10%  '

$ws.Range("E4").Value = '  +0.56%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.11'
$ws.Range("E5").Value = '  +0.46%  '

$ws.Range("E6").Value = '  -1.06%  '

$ws.Range("E7").Value = '  +0.60%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '26.84'
$ws.Range("E8").Value = '  +5.76%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '43.35'
$ws.Range("E9").Value = '  -5.34%  '

$ws.Range("E10").Value = '  +1.62%  '

$ws.Range("E11").Value = '  +1.19%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0908'
$ws.Range("E12").Value = '  +1.02%  '

$ws.Range("D13").Value = '1.805.51'
$ws.Range("E13").Value = '  +0.03%  '

$ws.Range("D14").Value = '1.595.06'
$ws.Range("E14").Value = '  +0.64%  '

$ws.Range("D15").Value = '29.546.27'
$ws.Range("E15").Value = '  +2.21%  '

$ws.Range("E16").Value = '  +3.08%  '

$ws.Range("E17").Value = '  +1.34%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '63.78'
$ws.Range("E18").Value = '  +2.78%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '240.41'
$ws.Range("E19").Value = '  +3.93%  '

$ws.Range("E20").Value = '  +2.31%  '

$ws.Range("E21").Value = '  +0.40%  '

$ws.Range("E22").Value = '  +0.81%  '

$ws.Range("E23").Value = '  +0.26%  '

$ws.Range("E24").Value = '  +1.28%  '

$ws.Range("E25").Value = '  +0.73%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '155.13'
$ws.Range("E26").Value = '  +1.78%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.33'
$ws.Range("E27").Value = '  +2.84%  '

$ws.Range("E28").Value = '  +0.70%  '

$ws.Range("E29").Value = '  +1.33%  '

$ws.Range("E30").Value = '  +0.53%  '

$ws.Range("E31").Value = '  +3.24%  '

$ws.Range("E32").Value = '  +0.64%  '

$ws.Range("E33").Value = '  +0.63%  '

$ws.Range("B34").Value = 'Maker'
$ws.Range("C34").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D34").Value = '1.436.39'
$ws.Range("E34").Value = '  +1.26%  '

$ws.Range("B35").Value = 'InternetComputer(DFINITY)'
$ws.Range("C35").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.14'
$ws.Range("E35").Value = '  +3.68%  '

$ws.Range("E36").Value = '  +3.06%  '

$ws.Range("E37").Value = '  -1.85%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.83'
$ws.Range("E38").Value = '  +3.00%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.30'
$ws.Range("E39").Value = '  +0.72%  '

$ws.Range("E40").Value = '  +1.52%  '

$ws.Range("E41").Value = '  +2.63%  '

$ws.Range("E42").Value = '  +0.89%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0491'
$ws.Range("E43").Value = '  +6.46%  '

$ws.Range("B44").Value = 'BitcoinSV'
$ws.Range("C44").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '53.52'
$ws.Range("E44").Value = '  +24.85%  '

$ws.Range("B45").Value = 'ARBITRUM'
$ws.Range("C45").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.800'
$ws.Range("E45").Value = '  +1.96%  '

$ws.Range("E46").Value = '  +0.58%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.978'
$ws.Range("E47").Value = '  +17.29%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '65.53'
$ws.Range("E48").Value = '  +1.87%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.32'
$ws.Range("E49").Value = '  +0.51%  '

$ws.Range("D50").Value = '1.737.12'
$ws.Range("E50").Value = '  +1.12%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '86.12'
$ws.Range("E51").Value = '  +1.09%  '
